$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.276194334030151
$ws.Range("B1").Value = 2.766923904418945
$ws.Range("C1").Value = 2.115238666534424
$ws.Range("D1").Value = 1.976908206939697
$ws.Range("E1").Value = 2.001348257064819
